# Reverse the order of the comma-separated "Recorded By" names in column G
# (e.g. "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com").
# Cells holding only a single name are left unchanged since reversing a
# single-element list is a no-op.

function Reverse-CommaList($s) {
    $parts = $s.Split(",")
    $n = $parts.Length
    if ($n -le 1) {
        return $s
    }
    $result = ""
    for ($i = $n - 1; $i -ge 0; $i--) {
        $p = $parts[$i].Trim()
        if ($result -ne "") {
            $result = $result + ", "
        }
        $result = $result + $p
    }
    return $result
}

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = 157
}

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $old = $cell.Text
    if ($old -ne "") {
        $new = Reverse-CommaList $old
        if ($new -ne $old) {
            $cell.Value = $new
        }
    }
}
